$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.768.25"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.591.00"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.77"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.348"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "3.046.65"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "60.784.33"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.63"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "2.598.83"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "352.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.425"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "2.712.56"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "0.0₃0839"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.35"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +11.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.31"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.75"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.18"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.941"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.87%  "
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.47"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.846"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.08"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0558"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.54"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0237"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.97"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.28%  "
